$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-9
$iValues = @(5, 8, 5, 7, 8, 5, 5, 1)
$jValues = @(6, 8, 5, 8, 8, 7, 6, 3)

for ($r = 0; $r -lt 8; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
